# "Fixed screen corruption in hospitals (MST collided with player select)"
#
# The "RAM windows cache" sheet lays out a chain of memory windows; each
# row's start offset (column B) is derived from another row's end offset
# (column C), and each row's size inputs live in column E (packed hex
# "XXYY" consumed by H/I/J). "MST in shop" (row 15) used to be placed
# right after "Shop items" (=C9), which put it on top of the "Player
# select" window -> corruption. Re-anchor it after "Player select"
# (=C22) instead, update its size input, and move "Yes/No" (row 23) off
# of the old anchor too so it no longer collides with the rest of the
# $DBF8 chain.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# MST in shop: re-anchor to start after "Player select" ends, and update its size input
# (leading apostrophe forces text, matching the column's quote-prefixed "look like
# a number but isn't" hex codes, e.g. 0320/032a would otherwise be read as numbers)
$ws.Range("E15").Value = "'0320"
$ws.Range("B15").Formula = "=C22"

# Enemy name window size input changed (row 12)
$ws.Range("E12").Value = "'032a"

# Yes/No: move off the old $DBF8-based anchor so it stops colliding there too
$ws.Range("B23").Formula = "=C22"

# Selection position as last left by the editor
$ws.Range("C14").Select()
